$d = $word.ActiveDocument

# --- 1. Remove the whole "Etendre longueur nav barre pied page" paragraph ---
$r1 = $d.Content
$r1.Find.Execute("Etendre longueur", $true, $false, $false, $false, $false,
                  $true, 1, $false, "", 0) | Out-Null
if ($r1.Find.Found) {
    $r1.Expand(4) | Out-Null
    $r1.Delete()
}

# --- 2. Mark the final tests as done: "TESTS finaux" -> "TESTS finaux FAIT" ---
$r2 = $d.Content
$r2.Find.Execute("TESTS finaux", $true, $false, $false, $false, $false,
                  $true, 1, $false, "", 0) | Out-Null
if ($r2.Find.Found) {
    $r2.Collapse(0)
    $r2.InsertAfter(" FAIT")
}

# --- 3. Move the "_GoBack" bookmark from the "Transfert de variable..." item
#        to the blank paragraph right after "TESTS finaux". ---
$targetIndex = -1
$i = 0
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    if ($p.Range.Text -like "TESTS finaux*") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -gt 0) {
    $blankAfterTests = $d.Paragraphs.Item($targetIndex + 1)

    if ($d.Bookmarks.Exists("_GoBack")) {
        $d.Bookmarks.Item("_GoBack").Delete()
    }
    $d.Bookmarks.Add("_GoBack", $blankAfterTests.Range)
}
